$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name and Link columns (plain text, order shifted)
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"

# Update Price and Volume(1h) columns - set as Text format to preserve exact string formatting
$deCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D38","E38","D39","E39","D40","E40","D41","E41","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $deCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "307.40"
$ws.Range("E2").Value = "0.00%"
$ws.Range("D3").Value = "41.44"
$ws.Range("E3").Value = "3.24%"
$ws.Range("D4").Value = "5.122"
$ws.Range("E4").Value = "2.34%"
$ws.Range("D5").Value = "0.07609"
$ws.Range("E5").Value = "-0.95%"
$ws.Range("D6").Value = "1.630"
$ws.Range("E6").Value = "0.74%"
$ws.Range("D7").Value = "2.475"
$ws.Range("E7").Value = "-2.88%"
$ws.Range("D8").Value = "0.9044"
$ws.Range("E8").Value = "1.23%"
$ws.Range("D9").Value = "0.1082"
$ws.Range("E9").Value = "9.47%"
$ws.Range("D10").Value = "0.1764"
$ws.Range("E10").Value = "1.26%"
$ws.Range("D11").Value = "0.09239"
$ws.Range("E11").Value = "3.85%"
$ws.Range("D12").Value = "0.04322"
$ws.Range("E12").Value = "-1.30%"
$ws.Range("D13").Value = "0.1052"
$ws.Range("E13").Value = "-0.41%"
$ws.Range("D14").Value = "0.001261"
$ws.Range("E14").Value = "-0.53%"
$ws.Range("D15").Value = "0.005803"
$ws.Range("E15").Value = "-1.93%"
$ws.Range("D16").Value = "3.362"
$ws.Range("E16").Value = "0.21%"
$ws.Range("D17").Value = "4.252"
$ws.Range("E17").Value = "-0.36%"
$ws.Range("D18").Value = "0.3296"
$ws.Range("E18").Value = "-1.98%"
$ws.Range("D19").Value = "6.576"
$ws.Range("E19").Value = "-6.44%"
$ws.Range("D20").Value = "0.1360"
$ws.Range("E20").Value = "0.79%"
$ws.Range("D21").Value = "0.2682"
$ws.Range("D22").Value = "0.04195"
$ws.Range("E22").Value = "-0.95%"
$ws.Range("D23").Value = "0.001221"
$ws.Range("E23").Value = "1.69%"
$ws.Range("D24").Value = "0.003994"
$ws.Range("E24").Value = "-1.81%"
$ws.Range("D25").Value = "0.0001298"
$ws.Range("E25").Value = "6.32%"
$ws.Range("D26").Value = "0.0003009"
$ws.Range("E26").Value = "0.92%"
$ws.Range("D38").Value = "0.02404"
$ws.Range("E38").Value = "2.46%"
$ws.Range("D39").Value = "0.05186"
$ws.Range("E39").Value = "0.26%"
$ws.Range("D40").Value = "0.007750"
$ws.Range("E40").Value = "-2.31%"
$ws.Range("D41").Value = "0.1299"
$ws.Range("E41").Value = "-1.63%"
$ws.Range("E42").Value = "6.20%"
$ws.Range("D43").Value = "0.001918"
$ws.Range("E43").Value = "-5.38%"
$ws.Range("D44").Value = "0.007870"
$ws.Range("E44").Value = "5.18%"
$ws.Range("D45").Value = "0.3049"
$ws.Range("E45").Value = "-8.09%"
$ws.Range("D46").Value = "0.00006754"
$ws.Range("E46").Value = "1.80%"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.21%"
$ws.Range("D48").Value = "0.004402"
$ws.Range("E48").Value = "-12.04%"
$ws.Range("D49").Value = "0.03223"
$ws.Range("E49").Value = "928.72%"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").Value = "-0.21%"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").Value = "-0.21%"
